$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null

# Row 75
$ws.Range("H75").Value = 32800
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 32800
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 32800
$ws.Range("N75").Value = -34672

# Row 78
$ws.Range("H78").Value = 32800
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 32800
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 98400
$ws.Range("N78").Value = -107760

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

# Row 103
$ws.Range("H103").Value = 9129.5
$ws.Range("I103").Value = 622.5
$ws.Range("J103").Value = 14800.833
$ws.Range("K103").Value = 1867.5
$ws.Range("L103").Value = 44402.499
$ws.Range("M103").Value = -1281.5
$ws.Range("N103").Value = -45574.499

# Row 111
$ws.Range("H111").Value = 1289.7142
$ws.Range("I111").Value = 1249.5
$ws.Range("J111").Value = 1305.8
$ws.Range("K111").Value = 3748.5
$ws.Range("L111").Value = 3917.4
$ws.Range("M111").Value = -681.5
$ws.Range("N111").Value = -10051.4

# Row 123
$ws.Range("H123").Value = 39000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 39000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -48800

# Row 129
$ws.Range("H129").Value = 850.9299999999999
$ws.Range("I129").Value = 347.75
$ws.Range("J129").Value = 894.68475
$ws.Range("K129").Value = 1043.25
$ws.Range("L129").Value = 2684.05425
$ws.Range("M129").Value = 3956.75
$ws.Range("N129").Value = -12684.05425

# Row 137
$ws.Range("H137").Value = 1192398.6
$ws.Range("I137").Value = 2507832.8
$ws.Range("J137").Value = 2244
$ws.Range("K137").Value = 7523498.399999999
$ws.Range("L137").Value = 6732
$ws.Range("M137").Value = -7520948.399999999

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
# Row 61
$ws.Range("H61").Value = 1229.9032
$ws.Range("I61").Value = 1160.6086
$ws.Range("J61").Value = 1429.125
$ws.Range("K61").Value = 1160.6086
$ws.Range("L61").Value = 1429.125
$ws.Range("M61").Value = -948.6086

# Row 74
$ws.Range("H74").Value = 5075.846
$ws.Range("I74").Value = 7705.0835
$ws.Range("J74").Value = 2822.2144
$ws.Range("K74").Value = 7705.0835
$ws.Range("L74").Value = 2822.2144
$ws.Range("M74").Value = -6831.0835

# Row 77
$ws.Range("H77").Value = 5075.846
$ws.Range("I77").Value = 7705.0835
$ws.Range("J77").Value = 2822.2144
$ws.Range("K77").Value = 38525.4175
$ws.Range("L77").Value = 14111.072
$ws.Range("M77").Value = -34157.4175

# Row 110
$ws.Range("H110").Value = 1475.7391
$ws.Range("I110").Value = 1366.1177
$ws.Range("J110").Value = 1786.3334
$ws.Range("K110").Value = 1366.1177
$ws.Range("L110").Value = 1786.3334
$ws.Range("M110").Value = 678.8823
$ws.Range("N110").Value = -5876.3334

# Row 136
$ws.Range("H136").Value = 1229.9032
$ws.Range("I136").Value = 1160.6086
$ws.Range("J136").Value = 1429.125
$ws.Range("K136").Value = 3481.8258
$ws.Range("L136").Value = 4287.375
$ws.Range("M136").Value = -931.8258000000001

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
# Row 7
$ws.Range("H7").Value = 1558982.8
$ws.Range("I7").Value = 5012499
$ws.Range("J7").Value = 24086.666
$ws.Range("K7").Value = 5012499
$ws.Range("L7").Value = 24086.666
$ws.Range("M7").Value = -5012386

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = $null

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = $null

# Row 103
$ws.Range("H103").Value = 38000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 38000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 38000
$ws.Range("N103").Value = -40344

# Row 107
$ws.Range("H107").Value = 1968.7142
$ws.Range("I107").Value = 1900.5625
$ws.Range("J107").Value = 2186.8
$ws.Range("K107").Value = 1900.5625
$ws.Range("L107").Value = 2186.8
$ws.Range("M107").Value = 19.4375
$ws.Range("N107").Value = -6026.8

# Row 118
$ws.Range("H118").Value = 28890
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 28890
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 28890
$ws.Range("N118").Value = -32204

# Row 134
$ws.Range("H134").Value = 2402.7908
$ws.Range("I134").Value = 1972.12
$ws.Range("J134").Value = 3000.9443
$ws.Range("K134").Value = 5916.36
$ws.Range("L134").Value = 9002.832900000001
$ws.Range("M134").Value = -3381.36
$ws.Range("N134").Value = -14072.8329

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 2651.6177
$ws.Range("I31").Value = 1143.4762
$ws.Range("J31").Value = 5087.846
$ws.Range("K31").Value = 1143.4762
$ws.Range("L31").Value = 5087.846
$ws.Range("M31").Value = -848.4762000000001

# Row 34
$ws.Range("H34").Value = 2651.6177
$ws.Range("I34").Value = 1143.4762
$ws.Range("J34").Value = 5087.846
$ws.Range("K34").Value = 1143.4762
$ws.Range("L34").Value = 5087.846
$ws.Range("M34").Value = -941.4762000000001

# Row 58
$ws.Range("H58").Value = 2400.738
$ws.Range("I58").Value = 1600.295
$ws.Range("J58").Value = 4523.6523
$ws.Range("K58").Value = 1600.295
$ws.Range("L58").Value = 4523.6523
$ws.Range("M58").Value = -1397.295
$ws.Range("N58").Value = -4929.6523

# Row 92
$ws.Range("H92").Value = 20500
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 20500
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 20500
$ws.Range("N92").Value = -25492

# Row 95
$ws.Range("H95").Value = 7500
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 7500
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 7500
$ws.Range("N95").Value = -12992

# Row 96
$ws.Range("H96").Value = 18206.666
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 18206.666
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 18206.666
$ws.Range("N96").Value = -23698.666

# Row 106
$ws.Range("H106").Value = 35398.332
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 35398.332
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 35398.332
$ws.Range("N106").Value = -37922.332

# Row 116
$ws.Range("H116").Value = 38500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 38500
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 38500
$ws.Range("N116").Value = -47678

# Row 122
$ws.Range("H122").Value = 1734.091
$ws.Range("I122").Value = 937.3077
$ws.Range("J122").Value = 2885
$ws.Range("K122").Value = 2811.9231
$ws.Range("L122").Value = 8655
$ws.Range("M122").Value = -361.9231
$ws.Range("N122").Value = -13555

# Row 132
$ws.Range("H132").Value = 2374.4644
$ws.Range("I132").Value = 1456.4762
$ws.Range("J132").Value = 5128.4287
$ws.Range("K132").Value = 4369.4286
$ws.Range("L132").Value = 15385.2861
$ws.Range("M132").Value = -1839.4286
$ws.Range("N132").Value = -20445.2861

# Row 134
$ws.Range("H134").Value = 5568.037
$ws.Range("I134").Value = 7818.067
$ws.Range("J134").Value = 2755.5
$ws.Range("K134").Value = 23454.201
$ws.Range("L134").Value = 8266.5
$ws.Range("M134").Value = -20919.201
$ws.Range("N134").Value = -13336.5

# Row 136
$ws.Range("H136").Value = 2400.738
$ws.Range("I136").Value = 1600.295
$ws.Range("J136").Value = 4523.6523
$ws.Range("K136").Value = 4800.885
$ws.Range("L136").Value = 13570.9569
$ws.Range("M136").Value = -2250.885
$ws.Range("N136").Value = -18670.9569

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
# Row 3
$ws.Range("H3").Value = 1739.4231
$ws.Range("I3").Value = 1484.375
$ws.Range("J3").Value = 4800
$ws.Range("K3").Value = 4453.125
$ws.Range("L3").Value = 14400
$ws.Range("M3").Value = -4341.125

# Row 9
$ws.Range("H9").Value = 155712.73
$ws.Range("I9").Value = 90250
$ws.Range("J9").Value = 234268
$ws.Range("K9").Value = 270750
$ws.Range("L9").Value = 702804
$ws.Range("M9").Value = -270526
$ws.Range("N9").Value = -703252

# Row 113
$ws.Range("H113").Value = 998.4
$ws.Range("I113").Value = 674.7273
$ws.Range("J113").Value = 1888.5
$ws.Range("K113").Value = 2024.1819
$ws.Range("L113").Value = 5665.5
$ws.Range("M113").Value = 145.8181
$ws.Range("N113").Value = -10005.5

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
# Row 132
$ws.Range("H132").Value = 2692.9167
$ws.Range("I132").Value = 1337.625
$ws.Range("J132").Value = 5403.5
$ws.Range("K132").Value = 4012.875
$ws.Range("L132").Value = 16210.5
$ws.Range("M132").Value = -1482.875

# Row 136
$ws.Range("H136").Value = 14854.714
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 14854.714
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 44564.142
$ws.Range("N136").Value = -49664.142

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
# Row 23
$ws.Range("H23").Value = 7252
$ws.Range("I23").Value = 5006
$ws.Range("J23").Value = 9498
$ws.Range("K23").Value = 5006
$ws.Range("L23").Value = 9498
$ws.Range("M23").Value = -4776
$ws.Range("N23").Value = -9958

# Row 54
$ws.Range("H54").Value = 30082.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 30082.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 30082.5
$ws.Range("N54").Value = -31370.5

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = $null

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = $null

# Row 136
$ws.Range("H136").Value = 4980.609
$ws.Range("I136").Value = 2421.4285
$ws.Range("J136").Value = 6100.25
$ws.Range("K136").Value = 7264.2855
$ws.Range("L136").Value = 18300.75
$ws.Range("M136").Value = -4714.2855
$ws.Range("N136").Value = -23400.75

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
# Row 82
$ws.Range("H82").Value = 42166.668
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 42166.668
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 42166.668
$ws.Range("N82").Value = -42932.668

# Row 85
$ws.Range("H85").Value = 42166.668
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 42166.668
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 42166.668
$ws.Range("N85").Value = -44818.668

# Row 92
$ws.Range("H92").Value = 25000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 25000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992

# Row 97
$ws.Range("H97").Value = 33300
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 33300
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 33300
$ws.Range("N97").Value = -35282

# Row 98
$ws.Range("H98").Value = 14935.4
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 14935.4
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 14935.4
$ws.Range("N98").Value = -20925.4

# Row 100
$ws.Range("H100").Value = 468.77777
$ws.Range("I100").Value = 464.875
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 929.75
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -388.75

# Row 135
$ws.Range("H135").Value = 32121.727
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 32121.727
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 32121.727
$ws.Range("N135").Value = -42261.727

# Row 136
$ws.Range("H136").Value = 5063.609
$ws.Range("I136").Value = 2859.875
$ws.Range("J136").Value = 10100.714
$ws.Range("K136").Value = 8579.625
$ws.Range("L136").Value = 30302.142
$ws.Range("M136").Value = -6029.625
$ws.Range("N136").Value = -35402.142
